# Indicator GGSTS217 ("familie of vrienden") is verwijderd uit het
# trendtabblad van het indicatoren overzicht, omdat er voor deze
# indicator geen trendcijfers beschikbaar zijn.

$wb = $excel.ActiveWorkbook

$wsTrends = $wb.Worksheets.Item("indicatoren trends")
$wsIndic  = $wb.Worksheets.Item("indicatoren")

# Row 24 holds indicator GGSTS217 - remove it; rows below shift up.
$wsTrends.Rows.Item(24).Delete()

# Rebuild the AutoFilter so its range shrinks from H35 to H34
# (toggling AutoFilterMode off first avoids merely switching the
# existing filter off).
$wsTrends.AutoFilterMode = $false
[void]$wsTrends.Range("A1:H34").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "indicatoren trends!_FilterDatabase") {
        $n.RefersTo = "='indicatoren trends'!`$A`$1:`$H`$34"
    }
}

# Make "indicatoren" the active / selected sheet instead of
# "indicatoren trends".
$wsIndic.Activate()
